$d = $word.ActiveDocument

function Find-ParaIndex($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $searchText) {
            return $i
        }
    }
    return -1
}

# Work from the bottom of the document upward so that earlier paragraph
# indices are never invalidated by edits made further down.

# 6) Insert a new empty paragraph right after the
#    "- Continuous refinement..." paragraph (near the end of the document,
#    before the final trailing empty paragraph).
$idxContinuous = Find-ParaIndex("- Continuous refinement and optimization of scheduling algorithms based on observed behavior and performance metrics are critical for developing robust and efficient solutions to the Unconstrained Examination Timetabling Problem.")
$d.Paragraphs.Item($idxContinuous).Range.InsertParagraphAfter()

# 5) Remove the "#### Overall Observations:" ... block, down through the
#    four trailing empty paragraphs right before "5. Shifts: 25, ...".
#    That heading paragraph through the last empty paragraph are deleted
#    outright, and the spacer paragraph ("   ") right before the heading
#    just has its text cleared (paragraph/run stays, text goes away).
$idxHeading = Find-ParaIndex("#### Overall Observations:")
$idxNextSection = Find-ParaIndex("5. Shifts: 25, Iterations: 300, Slots: 40")
$startDel = $d.Paragraphs.Item($idxHeading).Range.Start
$endDel = $d.Paragraphs.Item($idxNextSection - 1).Range.End
$d.Range($startDel, $endDel).Delete()

$idxSpacer3 = $idxHeading - 1
$d.Paragraphs.Item($idxSpacer3).Range.Text = ""

# 4) Clear the spacer paragraph ("   ") right after test case 2's bullets.
$idxTestCase3 = Find-ParaIndex("3. Shifts: 20, Iterations: 200, Slots: 30")
$idxSpacer2 = $idxTestCase3 - 1
$d.Paragraphs.Item($idxSpacer2).Range.Text = ""

# 3) Clear the spacer paragraph ("   ") right after test case 1's bullets.
$idxTestCase2 = Find-ParaIndex("2. Shifts: 5, Iterations: 50, Slots: 10")
$idxSpacer1 = $idxTestCase2 - 1
$d.Paragraphs.Item($idxSpacer1).Range.Text = ""

# 2) Fix the spacing in the "1. Shifts:10, Iterations: 100, Slots:20" line.
$d.Content.Find.Execute("1. Shifts:10, Iterations: 100, Slots:20", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1. Shifts: 10, Iterations: 100, Slots: 20", 2) | Out-Null

# 1) Remove the leading empty paragraph right before "### Observations:".
$idxObs = Find-ParaIndex("### Observations:")
$d.Paragraphs.Item($idxObs - 1).Range.Delete()
